$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3228075949727115
$ws.Range("C2").Value = 0.8267062650086889
$ws.Range("D2").Value = 1.14926112655861
$ws.Range("E2").Value = 1.072035972604749
$ws.Range("F2").Value = 1.050292850181011
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.2730085230503584
$ws.Range("C3").Value = 0.5944327505182431
$ws.Range("D3").Value = 0.6856664489529708
$ws.Range("E3").Value = 0.828049786518281
$ws.Range("F3").Value = 0.8044139377542835
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.1073812923613915
$ws.Range("C4").Value = 0.4856159529526489
$ws.Range("D4").Value = 0.418374132257862
$ws.Range("E4").Value = 0.6468184693233968
$ws.Range("F4").Value = 0.6574732710939286
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.2576619532183748
$ws.Range("C5").Value = 0.4703439316298302
$ws.Range("D5").Value = 0.3535200508351185
$ws.Range("E5").Value = 0.5945755215572859
$ws.Range("F5").Value = 0.5534188226639307
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.3604350276844117
$ws.Range("C6").Value = 0.3876405191120195
$ws.Range("D6").Value = 0.2123430664705153
$ws.Range("E6").Value = 0.460806973114031
$ws.Range("F6").Value = 0.297182586892516
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.3297347777357684
$ws.Range("C7").Value = 0.3594269622302418
$ws.Range("D7").Value = 0.198352598209372
$ws.Range("E7").Value = 0.4453679357670149
$ws.Range("F7").Value = 0.3106799050039986
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.362063388686015
$ws.Range("C8").Value = 0.3988911419790331
$ws.Range("D8").Value = 0.2152194845917584
$ws.Range("E8").Value = 0.4639175407243817
$ws.Range("F8").Value = 0.3018946605943391
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.3807292497673194
$ws.Range("C9").Value = 0.4290581283019066
$ws.Range("D9").Value = 0.2284558568240234
$ws.Range("E9").Value = 0.4779705606248396
$ws.Range("F9").Value = 0.3018146846158856
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.4249361656982711
$ws.Range("C10").Value = 0.4475969750223547
$ws.Range("D10").Value = 0.2365923687547577
$ws.Range("E10").Value = 0.4864076158478172
$ws.Range("F10").Value = 0.248241387000738
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.3961299661087644
$ws.Range("C11").Value = 0.4123773982225899
$ws.Range("D11").Value = 0.2050158272332239
$ws.Range("E11").Value = 0.4527867348247118
$ws.Range("F11").Value = 0.2311730404886565
$ws.Range("G11").Value = 10
